$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 801330
$ws.Range("A3").Value = 801835

$ws.Range("A4:A8").EntireRow.Delete()
